$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Estimated Effort Hours" values for the gantt chart fix
$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 55
$ws.Range("C5").Value = 60

# Update the selected cell on the sheet view
$ws.Range("C6").Select()
